# Generate Report for Handoff
#
# Updates the localization-status report: the original source file
# (2fdac74a-ee14-442c-a615-70e7b8763ed2.md) is replaced by a new source
# image (c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png) and two more files are
# now tracked: db8b0529-40d2-4137-86b8-89f951627e43.md (the original .md,
# now a dependency source) and ed541264-395a-45a1-b0f9-6b3230dddca9.png.

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------

# Row 2 - existing entry, filename + handoff date refreshed
$ovw.Range("A2").Hyperlinks.Delete()
$ovw.Hyperlinks.Add($ovw.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png", [Type]::Missing, [Type]::Missing, "c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png")
$ovw.Range("B2").Value = "Ready for handoff"
$ovw.Range("C2").Value = "Ready for handoff"
$ovw.Range("D2").Value = "2016-12-13 21:12:05"

# Row 3 - new entry: db8b0529-40d2-4137-86b8-89f951627e43.md
$ovw.Hyperlinks.Add($ovw.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/db8b0529-40d2-4137-86b8-89f951627e43.md", [Type]::Missing, [Type]::Missing, "db8b0529-40d2-4137-86b8-89f951627e43.md")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"
$ovw.Range("D3").Value = "2016-12-13 21:12:05"

# Row 4 - new entry: ed541264-395a-45a1-b0f9-6b3230dddca9.png
$ovw.Hyperlinks.Add($ovw.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/ed541264-395a-45a1-b0f9-6b3230dddca9.png", [Type]::Missing, [Type]::Missing, "ed541264-395a-45a1-b0f9-6b3230dddca9.png")
$ovw.Range("B4").Value = "Ready for handoff"
$ovw.Range("C4").Value = "Ready for handoff"
$ovw.Range("D4").Value = "2016-12-13 21:12:05"

# ---------------------------------------------------------------------
# zh-cn / de-de detail sheets share the same column layout:
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------

function Fill-DetailSheet($ws, $locale, $htBranch, $htHash, $handoffDate) {

    # Row 2 - c9460b57-...png (previously the .md source, same row)
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png", [Type]::Missing, [Type]::Missing, "c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png")

    $ws.Range("B2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/c9460b57-a1f1-428a-9bbf-eb9fa3ecdf3c.png", [Type]::Missing, [Type]::Missing, ".png")

    $ws.Range("C2").Value = "Ready for handoff"

    $ws.Range("D2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $htHash + "/ol-handoff/OpenLocalizationTestOrg/" + $htBranch + "/ci/ht/bd9228dec471118b35bba9445b64dbffea612527.png", [Type]::Missing, [Type]::Missing, "bd9228dec471118b35bba9445b64dbffea612527.png")

    $ws.Range("E2").Value = "2016-03-13 21:12:00"
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = "e2e\db8b0529-40d2-4137-86b8-89f951627e43.md"

    # Row 3 - db8b0529-...md (new dependency source)
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/db8b0529-40d2-4137-86b8-89f951627e43.md", [Type]::Missing, [Type]::Missing, "db8b0529-40d2-4137-86b8-89f951627e43.md")

    $ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/db8b0529-40d2-4137-86b8-89f951627e43.md", [Type]::Missing, [Type]::Missing, ".md")

    $ws.Range("C3").Value = "Ready for handoff"

    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $htHash + "/ol-handoff/OpenLocalizationTestOrg/" + $htBranch + "/ci/ht/db8b0529-40d2-4137-86b8-89f951627e43.07623363bb3cd05a5e56569a2a48e1bb4bec09f3." + $locale + ".xlf", [Type]::Missing, [Type]::Missing, "db8b0529-40d2-4137-86b8-89f951627e43.07623363bb3cd05a5e56569a2a48e1bb4bec09f3." + $locale + ".xlf")

    $ws.Range("E3").Value = $handoffDate
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    # Row 4 - ed541264-...png (new dependency)
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/ed541264-395a-45a1-b0f9-6b3230dddca9.png", [Type]::Missing, [Type]::Missing, "ed541264-395a-45a1-b0f9-6b3230dddca9.png")

    $ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9159691b2d6cf6544e2638afbe1f1fdd6cff1c3/e2e/ed541264-395a-45a1-b0f9-6b3230dddca9.png", [Type]::Missing, [Type]::Missing, ".png")

    $ws.Range("C4").Value = "Ready for handoff"

    $ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $htHash + "/ol-handoff/OpenLocalizationTestOrg/" + $htBranch + "/ci/ht/f403048fcf5f3da97bf9cde74252ccf97f38f5c4.png", [Type]::Missing, [Type]::Missing, "f403048fcf5f3da97bf9cde74252ccf97f38f5c4.png")

    $ws.Range("E4").Value = "2016-03-13 21:12:00"
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = "e2e\db8b0529-40d2-4137-86b8-89f951627e43.md"
}

Fill-DetailSheet $zhcn "zh-cn" "oltest.zh-cn" "38ecb33d90e69b50bf364eb4702dab7d43766f44" "2016-03-13 21:12:00"
Fill-DetailSheet $dede "de-de" "oltest.de-de" "9facae679b733400f7bb52e905f3c105c4891895" "2016-03-13 21:12:05"

# de-de sheet E column uses a slightly later handoff datetime than zh-cn
$dede.Range("E2").Value = "2016-03-13 21:12:05"
$dede.Range("E4").Value = "2016-03-13 21:12:05"
